$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1) to the new lower_case / underscore style names ---
$ws.Range("A1").Value = "colony_id"
$ws.Range("B1").Value = "diameter"
$ws.Range("C1").Value = "weight1"
$ws.Range("D1").Value = "wax_weight"
$ws.Range("E1").Value = "radius"
$ws.Range("F1").Value = "surface_area"
$ws.Range("G1").Value = "difference"

# --- Column widths (best-fit sizes picked up once the columns actually got used) ---
$ws.Columns.Item(5).ColumnWidth = 5.830729166666667
$ws.Columns.Item(6).ColumnWidth = 11.330729166666666
$ws.Columns.Item(7).ColumnWidth = 8.830729166666666
$ws.Columns.Item(15).ColumnWidth = 9.166666666666666
$ws.Columns.Item(16).ColumnWidth = 6.998697916666667
$ws.Columns.Item(17).ColumnWidth = 20.166666666666668

# --- Selection / scroll position: the sheet was last left with E15 selected ---
$ws.Range("E15").Select() | Out-Null
